$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.368.45'
$ws.Range("E2").Value = '  -1.71%  '

# Row 3
$ws.Range("D3").Value = '3.488.16'
$ws.Range("E3").Value = '  -1.99%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.86%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.43%  '

# Row 9
$ws.Range("E9").Value = '  -2.66%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.650'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.01%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.04'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.43%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000308'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.80%  '

# Row 13
$ws.Range("E13").Value = '  +0.83%  '

# Row 14
$ws.Range("D14").Value = '4.037.63'
$ws.Range("E14").Value = '  -2.16%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '600.39'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.26%  '

# Row 16
$ws.Range("D16").Value = '69.470.31'
$ws.Range("E16").Value = '  -1.75%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.58%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.85'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.60%  '

# Row 19
$ws.Range("D19").Value = '3.475.89'
$ws.Range("E19").Value = '  -2.69%  '

# Row 20
$ws.Range("E20").Value = '  -0.14%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.986'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.53%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.10'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.93%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '105.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +11.15%  '

# Row 24
$ws.Range("E24").Value = '  +1.95%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.90%  '

# Row 26
$ws.Range("E26").Value = '  +2.78%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.97%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.95%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.50%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.18%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.41'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.01%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +22.58%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.115'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.16%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.22'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.42%  '

# Row 35
$ws.Range("E35").Value = '  -6.79%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.997'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.27%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '524.58'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.60%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.396'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.17%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.64'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.10%  '

# Row 40
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '3.608.62'
$ws.Range("E40").Value = '  +0.21%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.70'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.91%  '

# Row 42
$ws.Range("D42").Value = '0.0₃0777'
$ws.Range("E42").Value = '  -2.72%  '

# Row 43
$ws.Range("E43").Value = '  +0.02%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0457'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.39%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.93'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.83%  '

# Row 46
$ws.Range("E46").Value = '  +3.04%  '

# Row 47
$ws.Range("E47").Value = '  -4.00%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.50%  '

# Row 49
$ws.Range("E49").Value = '  +0.32%  '

# Row 50
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000242'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.36%  '

# Row 51
$ws.Range("B51").Value = 'OceanProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.96%  '
